$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values would otherwise be auto-coerced to numbers by Excel
# are temporarily forced to Text format, then reset back to the Normal style
# so the saved cell carries no explicit style (matching the source data).

$ws.Range("D2").Value = "26.687.46"
$ws.Range("D3").Value = "1.595.20"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "1.818.35"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "1.598.19"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "26.662.05"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").Value = "1.298.39"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "1.731.16"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.904"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.37%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0985"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
